$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("templates")
$ws.Activate()

# Update the "estado" value for the crear_prod row (row 3) to "listo"
$ws.Range("E3").Value = "listo"

# Update the "estado3" value for the crear_prod row (row 3) to "listo"
$ws.Range("J3").Value = "listo"

# Update the active selection to E5
$ws.Range("E5").Select()
